# Load brand from database dynamic from production controller.
# User production repository extension to get all production brand.
#
# Updates the "Status" (K) and "Extended" (L) columns on Sheet1 with the
# latest task statuses / sample data, and moves the active selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Task statuses (column K)
$ws.Range("K5").Value  = "Done"
$ws.Range("K6").Value  = "Done"
$ws.Range("K7").Value  = "Working"
$ws.Range("K8").Value  = "Done"
$ws.Range("K9").Value  = "working"
$ws.Range("K11").Value = "Working"
$ws.Range("K14").Value = "Done"

# Sample data note (column L)
$ws.Range("L9").Value  = "Sample data"

# Move the active selection to match the author's final cursor position
$ws.Range("K12").Select()
